$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (date-like / numeric-like text) to remain plain text
$ws.Range("D17:D20").NumberFormat = "@"

# Row 17
$ws.Range("A17").Value = "UfYGAAAACAAJ"
$ws.Range("B17").Value = "El señor de los anillos"
$ws.Range("C17").Value = "Desconocido"
$ws.Range("D17").Value = "2002-02"
$ws.Range("E17").Value = "Kurt D. Bruner, Jim Ware"

# Row 18
$ws.Range("A18").Value = "WmdWtQAACAAJ"
$ws.Range("B18").Value = "El Señor de los anillos"
$ws.Range("C18").Value = "Desconocido"
$ws.Range("D18").Value = "2002"
$ws.Range("E18").Value = "J. R. R. Tolkien"

# Row 19
$ws.Range("A19").Value = "ZVwX0QEACAAJ"
$ws.Range("B19").Value = "El Señor de los Anillos"
$ws.Range("C19").Value = "Desconocido"
$ws.Range("D19").Value = "1985"
$ws.Range("E19").Value = "J. R. R. Tolkien"

# Row 20
$ws.Range("A20").Value = "ZcAlEAAAQBAJ"
$ws.Range("B20").Value = "The Lord of the Rings Illustrated"
$ws.Range("C20").Value = "A sumptuous slipcased edition of Tolkien's classic epic tale of adventure, fully illustrated in color by the author himself. This deluxe volume is quarterbound in leather and includes many special features unique to this edition. Since it was first published in 1954, The Lord of the Rings has been a book people have treasured. Steeped in unrivalled magic and otherworldliness, its sweeping fantasy and epic adventure has touched the hearts of young and old alike. Over 100 million copies of its many editions have been sold around the world, and occasional collectors' editions become prized and valuable items of publishing. This one-volume deluxe slipcased edition contains the complete text, fully corrected and reset, which is printed in red and black, and features thirty color illustrations, maps, and sketches drawn by Tolkien himself as he composed this epic work. These include the pages from the Book of Mazarbul, marvelous facsimiles created by Tolkien to accompany the famous `"Bridge of Khazad-dum`" chapter. Also appearing are two poster-size, fold-out maps revealing all the detail of Middle-earth. This very special deluxe edition is quarterbound in cloth and red leather, with raised ribs on the spine and stamped in two foils. The pages are edged in gold and contained within are special features unique to this edition."
$ws.Range("D20").Value = "2021-10-19"
$ws.Range("E20").Value = "J. R. R. Tolkien"
